# Update cryptocurrency price/volume data per the latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.834.36'
$ws.Range("E2").Value = '  -0.91%  '

# Row 3
$ws.Range("D3").Value = '1.634.69'
$ws.Range("E3").Value = '  -0.49%  '

# Row 4
$ws.Range("D4").Value = '''0.9991'
$ws.Range("E4").Value = '  -1.65%  '

# Row 5
$ws.Range("D5").Value = '''213.69'
$ws.Range("E5").Value = '  -1.26%  '

# Row 6
$ws.Range("D6").Value = '''0.5020'
$ws.Range("E6").Value = '  +0.47%  '

# Row 7
$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  -1.62%  '

# Row 8
$ws.Range("D8").Value = '''0.2555'
$ws.Range("E8").Value = '  -0.99%  '

# Row 9
$ws.Range("D9").Value = '''0.06356'
$ws.Range("E9").Value = '  -1.28%  '

# Row 10
$ws.Range("D10").Value = '''19.40'
$ws.Range("E10").Value = '  -0.57%  '

# Row 11
$ws.Range("D11").Value = '''0.07768'
$ws.Range("E11").Value = '  -0.15%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.646.01'
$ws.Range("E12").Value = '  +0.23%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.236'
$ws.Range("E13").Value = '  -0.52%  '

# Row 14
$ws.Range("D14").Value = '1.861.99'
$ws.Range("E14").Value = '  -0.19%  '

# Row 15
$ws.Range("D15").Value = '''0.5396'
$ws.Range("E15").Value = '  -1.33%  '

# Row 16
$ws.Range("D16").Value = '0.0₅7863'
$ws.Range("E16").Value = '  -1.20%  '

# Row 17
$ws.Range("D17").Value = '''64.29'
$ws.Range("E17").Value = '  +0.76%  '

# Row 18
$ws.Range("D18").Value = '25.873.53'
$ws.Range("E18").Value = '  -0.68%  '

# Row 19
$ws.Range("D19").Value = '''1.002'
$ws.Range("E19").Value = '  -1.46%  '

# Row 20
$ws.Range("D20").Value = '''195.41'
$ws.Range("E20").Value = '  -4.79%  '

# Row 21
$ws.Range("D21").Value = '''4.355'
$ws.Range("E21").Value = '  +0.87%  '

# Row 22
$ws.Range("D22").Value = '''9.872'
$ws.Range("E22").Value = '  -1.67%  '

# Row 23
$ws.Range("D23").Value = '''5.949'
$ws.Range("E23").Value = '  -0.52%  '

# Row 24
$ws.Range("D24").Value = '''1.004'
$ws.Range("E24").Value = '  -1.41%  '

# Row 25
$ws.Range("D25").Value = '''1.884'
$ws.Range("E25").Value = '  -4.73%  '

# Row 26
$ws.Range("D26").Value = '''139.33'
$ws.Range("E26").Value = '  -1.79%  '

# Row 27
$ws.Range("D27").Value = '''0.1133'
$ws.Range("E27").Value = '  -1.93%  '

# Row 28
$ws.Range("D28").Value = '''6.804'
$ws.Range("E28").Value = '  -0.23%  '

# Row 29
$ws.Range("D29").Value = '''15.61'
$ws.Range("E29").Value = '  -0.89%  '

# Row 30
$ws.Range("D30").Value = '''1.234'
$ws.Range("E30").Value = '  -0.97%  '

# Row 31
$ws.Range("D31").Value = '''0.04839'
$ws.Range("E31").Value = '  -4.10%  '

# Row 32
$ws.Range("D32").Value = '''3.237'
$ws.Range("E32").Value = '  -1.14%  '

# Row 33
$ws.Range("D33").Value = '''3.158'
$ws.Range("E33").Value = '  -1.64%  '

# Row 34
$ws.Range("D34").Value = '''1.524'
$ws.Range("E34").Value = '  -1.52%  '

# Row 35
$ws.Range("D35").Value = '''2.354'
$ws.Range("E35").Value = '  -0.35%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''0.8838'
$ws.Range("E36").Value = '  -1.22%  '

# Row 37
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '''2.596'
$ws.Range("E37").Value = '  -0.85%  '

# Row 38
$ws.Range("D38").Value = '''0.5490'
$ws.Range("E38").Value = '  -2.97%  '

# Row 39
$ws.Range("D39").Value = '1.122.62'
$ws.Range("E39").Value = '  -0.37%  '

# Row 40
$ws.Range("D40").Value = '''0.01556'
$ws.Range("E40").Value = '  -0.66%  '

# Row 41
$ws.Range("D41").Value = '''1.001'
$ws.Range("E41").Value = '  -1.80%  '

# Row 42
$ws.Range("D42").Value = '''5.648'
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$ws.Range("D43").Value = '''0.8105'
$ws.Range("E43").Value = '  -1.17%  '

# Row 44
$ws.Range("D44").Value = '''99.04'
$ws.Range("E44").Value = '  -0.93%  '

# Row 45
$ws.Range("E45").Value = '  +8.66%  '

# Row 46
$ws.Range("D46").Value = '1.773.22'
$ws.Range("E46").Value = '  -0.09%  '

# Row 47
$ws.Range("D47").Value = '''0.4505'
$ws.Range("E47").Value = '  -1.29%  '

# Row 48
$ws.Range("D48").Value = '''1.006'
$ws.Range("E48").Value = '  -0.84%  '

# Row 49
$ws.Range("D49").Value = '''54.86'
$ws.Range("E49").Value = '  -0.15%  '

# Row 50
$ws.Range("D50").Value = '''0.05033'
$ws.Range("E50").Value = '  -0.39%  '

# Row 51
$ws.Range("D51").Value = '''1.005'
$ws.Range("E51").Value = '  -0.96%  '
